# Generate Report for Handoff
#
# This updates the localization-status workbook so that the "b.md" file's
# handoff info is refreshed: status flips from "Handed back: in sync with
# en-US" to "Ready for handoff" on every sheet, the zh-cn/de-de detail
# sheets get a freshly generated handoff xliff name + timestamp for b.md,
# a "content duplicate" flag flips back to False, an explanatory error
# message is recorded, and the Error Detail column is widened to fit it.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- Overview sheet: row 3 is the b.md entry -----------------------------
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-17 22:36:20"

# --- zh-cn sheet ----------------------------------------------------------
# Row 2 = a.md, Row 3 = b.md (Source File Name in column A)
$zhcn.Range("C2").Value = "Ready for handoff"   # Status

$zhcn.Range("C3").Value = "Ready for handoff"   # Status
# Leading apostrophe forces text (otherwise "False" is parsed as a boolean);
# ClearFormats() drops the resulting quote-prefix formatting so the cell
# ends up a plain text cell, same as the untouched "False" cells elsewhere.
$zhcn.Range("F3").Value = "'False"              # Content Duplicate
$zhcn.Range("F3").ClearFormats()
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"  # Latest Handoff File
$zhcn.Range("H3").Value = "2016-08-17 22:36:13" # Latest Handoff Datetime
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8a1b7a54e682ba751164359b31bf1281f8d08ffd/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f48943ffc4adf18eb3ef2d9c67a7c944412d107b/e2e/b.md."

# Widen the Error Detail column (P / col 16) to fit the new message.
$zhcn.Columns.Item(16).ColumnWidth = 39.15

# --- de-de sheet -----------------------------------------------------------
$dede.Range("C2").Value = "Ready for handoff"   # Status

$dede.Range("C3").Value = "Ready for handoff"   # Status
$dede.Range("F3").Value = "'False"              # Content Duplicate
$dede.Range("F3").ClearFormats()
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"  # Latest Handoff File
$dede.Range("H3").Value = "2016-08-17 22:36:20" # Latest Handoff Datetime
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8a1b7a54e682ba751164359b31bf1281f8d08ffd/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f48943ffc4adf18eb3ef2d9c67a7c944412d107b/e2e/b.md."

$dede.Columns.Item(16).ColumnWidth = 39.15
